# Update Data Sources from LFX — refresh the table style applied to every
# data-source table in the deck.
#
# The six "Data Sources" tables (originally generated by the Google Slides ->
# PPTX converter) all reference the old custom table style
# {50C352CE-1B1D-4D2E-B4E8-68B4B13091BA}. The automated LFX data refresh
# re-applies the (now current) table style
# {27CCD006-CB68-44EF-8F73-6356A743C068} to each of those tables, leaving
# everything else (grid, rows, cell text, fills, ...) untouched.

$OldStyleId = "{50C352CE-1B1D-4D2E-B4E8-68B4B13091BA}"
$NewStyleId = "{27CCD006-CB68-44EF-8F73-6356A743C068}"

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style.Name -eq $OldStyleId) {
                $table.ApplyStyle($NewStyleId, $true)
            }
        }
    }
}
